$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows at position 620, shifting existing rows 620-710
# down to 624-714 (mirrors the diff: dimension A1:T710 -> A1:T714).
$ws.Range("620:623").Insert()

# Common (constant across all data rows in this sheet) column values.
$mercadoId = 6
$mercado   = "Mercado Mayorista Lo Valledor de Santiago"
$region    = "Metropolitana"
$codreg    = 13
$tipo      = "Fruta"
$productoId = 100108
$producto   = "Tropicales y subtropicales"
$categoriaId = 100108005
$categoria   = "Piña"
$variedad    = "Caramelo"
$origen      = "Ecuador"

function Set-PinaRow {
    param($row, $fecha, $calidad, $volumen, $precioMin, $precioMax, $precioProm, $unidad, $precioKg, $kgUnidad)

    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $tipo
    $ws.Cells.Item($row, 7).Value  = $productoId
    $ws.Cells.Item($row, 8).Value  = $producto
    $ws.Cells.Item($row, 9).Value  = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $calidad
    $ws.Cells.Item($row, 13).Value = $volumen
    $ws.Cells.Item($row, 14).Value = $precioMin
    $ws.Cells.Item($row, 15).Value = $precioMax
    $ws.Cells.Item($row, 16).Value = $precioProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $precioKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}

Set-PinaRow 620 44474 "Especial" 461 18000 20000 18863 "`$/caja 10 unidades" 1886 10
Set-PinaRow 621 44474 "Primera"  471 18000 20000 18463 "`$/caja 12 unidades" 1539 12
Set-PinaRow 622 44474 "Segunda"  363 18000 20000 18452 "`$/caja 14 unidades" 1318 14
Set-PinaRow 623 44474 "Tercera"  35  20000 20000 20000 "`$/caja 16 unidades" 1250 16
